# Insert a new weekly price record for "Vega Monumental Concepción - Pepino
# ensalada" as row 144, pushing the previously existing rows 144-174 down to
# 145-175 (dimension grows from A1:R174 to A1:R175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 144:174 down by inserting a blank row above row 144.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new record's data.
$ws.Cells.Item(144, 1).Value = 11
$ws.Cells.Item(144, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(144, 3).Value = "Bíobío"
$ws.Cells.Item(144, 4).Value = 44932
$ws.Cells.Item(144, 5).Value = 8
$ws.Cells.Item(144, 6).Value = 100112043
$ws.Cells.Item(144, 7).Value = "Pepino ensalada"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 2000
$ws.Cells.Item(144, 11).Value = 900
$ws.Cells.Item(144, 12).Value = 1000
$ws.Cells.Item(144, 13).Value = 950
$ws.Cells.Item(144, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(144, 15).Value = "Región Metropolitana"
$ws.Cells.Item(144, 16).Value = 16
$ws.Cells.Item(144, 17).Value = 60
$ws.Cells.Item(144, 18).Value = "Hortaliza"
